$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-NumericLookingText($range, [string]$text) {
    # Column C holds purely-numeric-looking labels ("2"). A plain
    # $range.Value = "2" gets auto-coerced to a Number by the COM layer,
    # but the source data is a text label (inlineStr "2"), so force text
    # the way Excel's own leading-apostrophe input does, then clear the
    # resulting quote-prefix style back to Normal so no per-cell style
    # attribute is left behind.
    $range.Value = "'" + $text
    $range.Style = "Normal"
}

# Update existing rows 2 and 3 (ImagesCount column C) from 4 to 2
Set-NumericLookingText $ws.Range("C2") "2"
Set-NumericLookingText $ws.Range("C3") "2"

# Add new rows 4-10 (w3..w9), mirroring the existing w1/w2 pattern
$data = @(
    @("w3", "5-6", "2"),
    @("w4", "7-8", "2"),
    @("w5", "9-10", "2"),
    @("w6", "11-12", "2"),
    @("w7", "13-14", "2"),
    @("w8", "15-16", "2"),
    @("w9", "17-18", "2")
)

$row = 4
foreach ($entry in $data) {
    # Columns A/B ("w3", "5-6", ...) are not numeric-looking, so a plain
    # assignment already round-trips as text with no style side effects.
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    Set-NumericLookingText $ws.Cells.Item($row, 3) $entry[2]
    $row++
}

# Update selection to D8
$ws.Range("D8").Select()

# Match the author's resized Excel window (best effort - harmless if the
# host doesn't persist window geometry into bookViews/workbookView).
try {
    $w = $excel.ActiveWindow
    $w.Left = 1820
    $w.Top = 1820
    $w.Width = 14400
    $w.Height = 7360
} catch {
}
